$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the TC_ID value from column A for rows 4 through 11 (moved to outside
# the repository directory - the repeated TC_ID label is no longer needed on
# these rows).
$ws.Range("A4:A11").Clear()
